$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph and before the "Gameplay Mechanics" Heading2 paragraph. The
#    injected WordOpenXML fragment carries no <w:pPr>, so it overrides
#    whatever heading style the blank placeholder paragraph would otherwise
#    have inherited, leaving a plain/body paragraph exactly like the diff.
# ---------------------------------------------------------------------------

$titlePara0 = $d.Paragraphs.Item(1)   # "Play Cosmic Heart Free: ..." (Heading1)
$titlePara0.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)     # new blank paragraph, now sits before "Gameplay Mechanics"
$metaRange = $metaPara.Range
$metaRange.Collapse(1)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the universe and win big with Cosmic Heart, a high RTP, high volatility slot game with innovative special features. Try this game for free now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRange.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold title paragraph
#    ("Play Cosmic Heart Free: ...") and rewrite the italic paragraph that
#    used to hold the meta-description text so that it now holds the image
#    prompt text instead.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)          # italic meta-description paragraph
$titlePara = $d.Paragraphs.Item($count - 1)     # bold duplicated title paragraph

# Insert the replacement (image-prompt) paragraph right before the old
# italic paragraph, so it inherits the same (style-less) paragraph
# formatting, then delete the two paragraphs that are no longer wanted.
$lastPara.Range.InsertParagraphBefore() | Out-Null

$imgPara = $d.Paragraphs.Item($count)
$imgRange = $imgPara.Range
$imgRange.Collapse(1)

$imgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Create a cartoon style feature image for the slot game &quot;Cosmic Heart&quot;. The image should feature a happy Maya warrior with glasses. The Maya warrior could be seen in a spaceship or on a planet, surrounded by elements of outer space such as stars, planets, or galaxies. The image should be colorful, eye-catching, and highlight the theme of space adventure and exploration.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$imgRange.InsertXML($imgXml) | Out-Null

# The document now looks like:
#   ... [count-1] bold "Play Cosmic Heart Free..." title paragraph
#       [count]   new italic image-prompt paragraph
#       [count+1] old italic "Explore the universe..." paragraph (now stale)
# Delete the stale trailing paragraph first (keeps indices simple), then the
# duplicated bold title paragraph.

$staleOldItalic = $d.Paragraphs.Item($count + 1)
$staleOldItalic.Range.Delete() | Out-Null

$dupTitle = $d.Paragraphs.Item($count - 1)
$dupTitle.Range.Delete() | Out-Null

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
